# "build the sv events from vcf files"
#
# The patient-ID column (A) used to hold the raw VCF filenames
# (e.g. "LP6005878-DNA_A01.SV.vcf.gz"). Strip the ".SV.vcf.gz" suffix so the
# column holds plain patient IDs (e.g. "LP6005878-DNA_A01") that can be used
# to key the per-sample SV event data pulled from those VCFs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $raw = $cell.Value2
    if ($raw -ne $null) {
        $cell.Value = $raw.Replace(".SV.vcf.gz", "")
    }
}

# Widen the patient-ID column now that it is the meaningful lookup key, and
# reset the view back to the top of the sheet.
$ws.Columns.Item(1).ColumnWidth = 25.66666666666667

$ws.Range("A2").Select()
